$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24).
$ws.Rows.Item(13).Insert()

# Row 13 (new): only B/C filled with the teacher name, no label in A.
# Copy number/text formatting from row 14's B/C cells (style already correct: s=2/s=3)
# onto the brand-new row 13 cells before setting their values, so they do not
# inherit the default bold style.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840917 - Fabricio Maciel Gomes"
$ws.Range("C13").Value = "5840917 - Fabricio Maciel Gomes"

# Row 10 (Objetivos:) - replace the objective text
$ws.Range("B10").Value = "Proporcionar conhecimento de Pesquisa Operacional como ciência aplicada."
$ws.Range("C10").Value = "Proporcionar conhecimento de Pesquisa Operacional como ciência aplicada."

# Row 14 (Programa resumido:) - replace "Semestral" with the real summary
$ws.Range("B14").Value = "Modelos PERT/COM, Programação Linear Inteira, Programação Dinâmica, Métodos Heurísticos, Modelos e Técnicas de Previsão."
$ws.Range("C14").Value = "Modelos PERT/COM, Programação Linear Inteira, Programação Dinâmica, Métodos Heurísticos, Modelos e Técnicas de Previsão."

# Row 16 (Programa:) - replace with full Portuguese syllabus text
$ws.Range("B16").Value = "1. Modelos PERT/COM2. Programação Linear Inteira; 2.1. Modelamento de problemas de PLI. 2.2 Algoritmo de ramificação e avaliação progressiva (branchand-bound).3. Programação Dinâmica3. Métodos Heurísticos; 3.1. Algoritmos Genéticos; 3.2 Recozimento Simulado; 3.3 Aplicação em problemas de otimização.4. Modelos e Técnicas de Previsão"
$ws.Range("C16").Value = "1. Modelos PERT/COM2. Programação Linear Inteira; 2.1. Modelamento de problemas de PLI. 2.2 Algoritmo de ramificação e avaliação progressiva (branchand-bound).3. Programação Dinâmica3. Métodos Heurísticos; 3.1. Algoritmos Genéticos; 3.2 Recozimento Simulado; 3.3 Aplicação em problemas de otimização.4. Modelos e Técnicas de Previsão"

# Row 19 (Método:) - now holds the evaluation method text (was shifted from old row19's Critério text)
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 (Critério:) - now holds the "NF>=5,0." text
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21 (Norma de recuperação:) - now holds the "Média aritmética..." text
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."

# Row 22 (Bibliografia:) - now holds the actual bibliography text
$ws.Range("B22").Value = "1. HILLIER, F.S., LIEBERMAN, G.J., “Introdução à Pesquisa Operacional”, 8ªed., Editora McGraw-Hill, 2006.2. LACHTERMACHER, G., “Pesquisa Operacional na Tomada de Decisão (modelagem em Excel)”, 4ª ed., Editora Campus, 2009.3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., “An Introduction to Management Science” 9ª ed., South-Western College Publishing, 2000.4. PIZZOLATO, N. D. e GANDOLPHO, A. A. “Técnicas de Otimização”, LTC Editora, 2009.5. TAHA, H. A ., “Pesquisa Operacional”, 8ª ed., Pearson/Prentice Hall, 2008."
$ws.Range("C22").Value = "1. HILLIER, F.S., LIEBERMAN, G.J., “Introdução à Pesquisa Operacional”, 8ªed., Editora McGraw-Hill, 2006.2. LACHTERMACHER, G., “Pesquisa Operacional na Tomada de Decisão (modelagem em Excel)”, 4ª ed., Editora Campus, 2009.3. ANDERSON, D.R., SWEENEY, D.J. e WILLIAMS, T.A., “An Introduction to Management Science” 9ª ed., South-Western College Publishing, 2000.4. PIZZOLATO, N. D. e GANDOLPHO, A. A. “Técnicas de Otimização”, LTC Editora, 2009.5. TAHA, H. A ., “Pesquisa Operacional”, 8ª ed., Pearson/Prentice Hall, 2008."
